$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (shifts RL1 and everything below down by one,
# and Excel auto-updates the SUM formulas in the totals row).
$ws.Rows("5:5").Insert()

# Fill in the new row 5 with the capacitor (C1) added across the strobe relay power.
$ws.Range("A5").Value = "C1"
$ws.Range("B5").Value = "Ceramic Capacitor"
$ws.Range("C5").Value = "0.47nF; 50v"
$ws.Range("D5").Value = "0805B471K500CT"
$ws.Range("F5").Value = "SMD"
$ws.Range("G5").Value = "'0805"
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0.1

# Grow the BOM table so the new row is included within it.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J10"))
